$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.326.55'
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = '1.622.30'
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.49'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +1.11%  '
$ws.Range("E8").Value = '  +1.58%  '
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.91'
$ws.Range("E10").Value = '  +5.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0816'
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '1.847.81'
$ws.Range("D13").Value = '1.627.67'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").Value = '26.334.88'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.50'
$ws.Range("E17").Value = '  +4.16%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.74'
$ws.Range("E20").Value = '  +1.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.29'
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("E24").Value = '  +6.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.34'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.20'
$ws.Range("E28").Value = '  +0.93%  '
$ws.Range("E29").Value = '  +2.02%  '
$ws.Range("E30").Value = '  +11.33%  '
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("E32").Value = '  +2.84%  '
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("E34").Value = '  +1.96%  '
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("D36").Value = '1.180.99'
$ws.Range("E36").Value = '  +5.15%  '
$ws.Range("E37").Value = '  +1.76%  '
$ws.Range("E38").Value = '  +3.65%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.495'
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.789'
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.34'
$ws.Range("E43").Value = '  +5.14%  '
$ws.Range("D44").Value = '1.759.32'
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.64'
$ws.Range("E45").Value = '  +1.29%  '
$ws.Range("E46").Value = '  +15.83%  '
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.00'
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("E51").Value = '  -0.40%  '
